$d = $word.ActiveDocument

$d.Content.Find.Execute("98-13=85", $true, $true, $false, $false, $false, $true, 1, $false, "25+21=46", 2) | Out-Null
$d.Content.Find.Execute("38+49=87", $true, $true, $false, $false, $false, $true, 1, $false, "22-4=18", 2) | Out-Null
$d.Content.Find.Execute("1+81=82", $true, $true, $false, $false, $false, $true, 1, $false, "98-85=13", 2) | Out-Null
$d.Content.Find.Execute("95-64=31", $true, $true, $false, $false, $false, $true, 1, $false, "87-1=86", 2) | Out-Null
$d.Content.Find.Execute("81-70=11", $true, $true, $false, $false, $false, $true, 1, $false, "42+11=53", 2) | Out-Null
$d.Content.Find.Execute("8+59=67", $true, $true, $false, $false, $false, $true, 1, $false, "61+35=96", 2) | Out-Null
$d.Content.Find.Execute("39+37=76", $true, $true, $false, $false, $false, $true, 1, $false, "52+13=65", 2) | Out-Null
$d.Content.Find.Execute("77-14=63", $true, $true, $false, $false, $false, $true, 1, $false, "13-10=3", 2) | Out-Null
$d.Content.Find.Execute("32+63=95", $true, $true, $false, $false, $false, $true, 1, $false, "77-23=54", 2) | Out-Null
$d.Content.Find.Execute("59+23=82", $true, $true, $false, $false, $false, $true, 1, $false, "44+36=80", 2) | Out-Null
$d.Content.Find.Execute("43-42=1", $true, $true, $false, $false, $false, $true, 1, $false, "85-75=10", 2) | Out-Null
$d.Content.Find.Execute("55-29=26", $true, $true, $false, $false, $false, $true, 1, $false, "36+59=95", 2) | Out-Null
$d.Content.Find.Execute("96+0=96", $true, $true, $false, $false, $false, $true, 1, $false, "68-66=2", 2) | Out-Null
$d.Content.Find.Execute("8+1=9", $true, $true, $false, $false, $false, $true, 1, $false, "43+56=99", 2) | Out-Null
$d.Content.Find.Execute("81-29=52", $true, $true, $false, $false, $false, $true, 1, $false, "14+26=40", 2) | Out-Null
$d.Content.Find.Execute("63-26=37", $true, $true, $false, $false, $false, $true, 1, $false, "59-40=19", 2) | Out-Null
$d.Content.Find.Execute("37+10=47", $true, $true, $false, $false, $false, $true, 1, $false, "82-23=59", 2) | Out-Null
$d.Content.Find.Execute("69-38=31", $true, $true, $false, $false, $false, $true, 1, $false, "93-13=80", 2) | Out-Null
$d.Content.Find.Execute("96-93=3", $true, $true, $false, $false, $false, $true, 1, $false, "66-56=10", 2) | Out-Null
$d.Content.Find.Execute("27+40=67", $true, $true, $false, $false, $false, $true, 1, $false, "83+7=90", 2) | Out-Null
$d.Content.Find.Execute("81+16=97", $true, $true, $false, $false, $false, $true, 1, $false, "83-10=73", 2) | Out-Null
$d.Content.Find.Execute("42+53=95", $true, $true, $false, $false, $false, $true, 1, $false, "88+3=91", 2) | Out-Null
$d.Content.Find.Execute("32-24=8", $true, $true, $false, $false, $false, $true, 1, $false, "95-35=60", 2) | Out-Null
$d.Content.Find.Execute("56+35=91", $true, $true, $false, $false, $false, $true, 1, $false, "54+21=75", 2) | Out-Null
$d.Content.Find.Execute("23+14=37", $true, $true, $false, $false, $false, $true, 1, $false, "50+41=91", 2) | Out-Null
$d.Content.Find.Execute("39+14=53", $true, $true, $false, $false, $false, $true, 1, $false, "37+6=43", 2) | Out-Null
$d.Content.Find.Execute("0+50=50", $true, $true, $false, $false, $false, $true, 1, $false, "74-33=41", 2) | Out-Null
$d.Content.Find.Execute("76+2=78", $true, $true, $false, $false, $false, $true, 1, $false, "26+64=90", 2) | Out-Null
$d.Content.Find.Execute("68+21=89", $true, $true, $false, $false, $false, $true, 1, $false, "49+17=66", 2) | Out-Null
$d.Content.Find.Execute("69-24=45", $true, $true, $false, $false, $false, $true, 1, $false, "18+38=56", 2) | Out-Null
$d.Content.Find.Execute("12+68=80", $true, $true, $false, $false, $false, $true, 1, $false, "88-33=55", 2) | Out-Null
$d.Content.Find.Execute("2+42=44", $true, $true, $false, $false, $false, $true, 1, $false, "9+19=28", 2) | Out-Null
$d.Content.Find.Execute("89-10=79", $true, $true, $false, $false, $false, $true, 1, $false, "12+62=74", 2) | Out-Null
$d.Content.Find.Execute("42+57=99", $true, $true, $false, $false, $false, $true, 1, $false, "37-24=13", 2) | Out-Null
$d.Content.Find.Execute("78+0=78", $true, $true, $false, $false, $false, $true, 1, $false, "47-47=0", 2) | Out-Null
$d.Content.Find.Execute("14+12=26", $true, $true, $false, $false, $false, $true, 1, $false, "23+36=59", 2) | Out-Null
$d.Content.Find.Execute("68-34=34", $true, $true, $false, $false, $false, $true, 1, $false, "4+16=20", 2) | Out-Null
$d.Content.Find.Execute("88-49=39", $true, $true, $false, $false, $false, $true, 1, $false, "52-21=31", 2) | Out-Null
$d.Content.Find.Execute("19+24=43", $true, $true, $false, $false, $false, $true, 1, $false, "54-43=11", 2) | Out-Null
$d.Content.Find.Execute("89-5=84", $true, $true, $false, $false, $false, $true, 1, $false, "27-25=2", 2) | Out-Null
$d.Content.Find.Execute("58+27=85", $true, $true, $false, $false, $false, $true, 1, $false, "66+23=89", 2) | Out-Null
$d.Content.Find.Execute("39+36=75", $true, $true, $false, $false, $false, $true, 1, $false, "17+67=84", 2) | Out-Null
$d.Content.Find.Execute("11+77=88", $true, $true, $false, $false, $false, $true, 1, $false, "91-11=80", 2) | Out-Null
$d.Content.Find.Execute("4+39=43", $true, $true, $false, $false, $false, $true, 1, $false, "5+53=58", 2) | Out-Null
$d.Content.Find.Execute("65-45=20", $true, $true, $false, $false, $false, $true, 1, $false, "2+67=69", 2) | Out-Null
$d.Content.Find.Execute("98-88=10", $true, $true, $false, $false, $false, $true, 1, $false, "78-66=12", 2) | Out-Null
$d.Content.Find.Execute("86-70=16", $true, $true, $false, $false, $false, $true, 1, $false, "17+66=83", 2) | Out-Null
$d.Content.Find.Execute("19+3=22", $true, $true, $false, $false, $false, $true, 1, $false, "27+61=88", 2) | Out-Null
$d.Content.Find.Execute("39+32=71", $true, $true, $false, $false, $false, $true, 1, $false, "81-54=27", 2) | Out-Null
$d.Content.Find.Execute("97-11=86", $true, $true, $false, $false, $false, $true, 1, $false, "48+47=95", 2) | Out-Null
$d.Content.Find.Execute("80-2=78", $true, $true, $false, $false, $false, $true, 1, $false, "45-2=43", 2) | Out-Null
$d.Content.Find.Execute("96-72=24", $true, $true, $false, $false, $false, $true, 1, $false, "61-17=44", 2) | Out-Null
$d.Content.Find.Execute("57-4=53", $true, $true, $false, $false, $false, $true, 1, $false, "27+20=47", 2) | Out-Null
$d.Content.Find.Execute("63-19=44", $true, $true, $false, $false, $false, $true, 1, $false, "19+1=20", 2) | Out-Null
$d.Content.Find.Execute("94-17=77", $true, $true, $false, $false, $false, $true, 1, $false, "63+5=68", 2) | Out-Null
$d.Content.Find.Execute("41+39=80", $true, $true, $false, $false, $false, $true, 1, $false, "95-74=21", 2) | Out-Null
$d.Content.Find.Execute("54+19=73", $true, $true, $false, $false, $false, $true, 1, $false, "42+47=89", 2) | Out-Null
$d.Content.Find.Execute("31+7=38", $true, $true, $false, $false, $false, $true, 1, $false, "24+61=85", 2) | Out-Null
$d.Content.Find.Execute("69-35=34", $true, $true, $false, $false, $false, $true, 1, $false, "69-33=36", 2) | Out-Null
$d.Content.Find.Execute("87+3=90", $true, $true, $false, $false, $false, $true, 1, $false, "75-71=4", 2) | Out-Null
$d.Content.Find.Execute("22-12=10", $true, $true, $false, $false, $false, $true, 1, $false, "17+11=28", 2) | Out-Null
$d.Content.Find.Execute("23+13=36", $true, $true, $false, $false, $false, $true, 1, $false, "76-30=46", 2) | Out-Null
$d.Content.Find.Execute("92-33=59", $true, $true, $false, $false, $false, $true, 1, $false, "83-75=8", 2) | Out-Null
$d.Content.Find.Execute("34-2=32", $true, $true, $false, $false, $false, $true, 1, $false, "54+38=92", 2) | Out-Null
$d.Content.Find.Execute("61+19=80", $true, $true, $false, $false, $false, $true, 1, $false, "93-40=53", 2) | Out-Null
$d.Content.Find.Execute("53-22=31", $true, $true, $false, $false, $false, $true, 1, $false, "30+17=47", 2) | Out-Null
$d.Content.Find.Execute("0+24=24", $true, $true, $false, $false, $false, $true, 1, $false, "73-4=69", 2) | Out-Null
$d.Content.Find.Execute("67-57=10", $true, $true, $false, $false, $false, $true, 1, $false, "37-1=36", 2) | Out-Null
$d.Content.Find.Execute("24+35=59", $true, $true, $false, $false, $false, $true, 1, $false, "36-33=3", 2) | Out-Null
$d.Content.Find.Execute("48+26=74", $true, $true, $false, $false, $false, $true, 1, $false, "2+17=19", 2) | Out-Null
$d.Content.Find.Execute("64-59=5", $true, $true, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$d.Content.Find.Execute("16-0=16", $true, $true, $false, $false, $false, $true, 1, $false, "8+81=89", 2) | Out-Null
$d.Content.Find.Execute("37+57=94", $true, $true, $false, $false, $false, $true, 1, $false, "8+19=27", 2) | Out-Null
$d.Content.Find.Execute("10-4=6", $true, $true, $false, $false, $false, $true, 1, $false, "5+79=84", 2) | Out-Null
$d.Content.Find.Execute("68-40=28", $true, $true, $false, $false, $false, $true, 1, $false, "2+17=19", 2) | Out-Null
$d.Content.Find.Execute("57+15=72", $true, $true, $false, $false, $false, $true, 1, $false, "66-9=57", 2) | Out-Null
$d.Content.Find.Execute("38+44=82", $true, $true, $false, $false, $false, $true, 1, $false, "79-3=76", 2) | Out-Null
$d.Content.Find.Execute("45-41=4", $true, $true, $false, $false, $false, $true, 1, $false, "6+37=43", 2) | Out-Null
$d.Content.Find.Execute("85+12=97", $true, $true, $false, $false, $false, $true, 1, $false, "25+67=92", 2) | Out-Null
$d.Content.Find.Execute("94-26=68", $true, $true, $false, $false, $false, $true, 1, $false, "40-8=32", 2) | Out-Null
$d.Content.Find.Execute("42-6=36", $true, $true, $false, $false, $false, $true, 1, $false, "7+11=18", 2) | Out-Null
$d.Content.Find.Execute("19+73=92", $true, $true, $false, $false, $false, $true, 1, $false, "12+26=38", 2) | Out-Null
$d.Content.Find.Execute("13+69=82", $true, $true, $false, $false, $false, $true, 1, $false, "58+33=91", 2) | Out-Null
$d.Content.Find.Execute("8+36=44", $true, $true, $false, $false, $false, $true, 1, $false, "10-3=7", 2) | Out-Null
$d.Content.Find.Execute("99-28=71", $true, $true, $false, $false, $false, $true, 1, $false, "63+29=92", 2) | Out-Null
$d.Content.Find.Execute("31+14=45", $true, $true, $false, $false, $false, $true, 1, $false, "6+79=85", 2) | Out-Null
$d.Content.Find.Execute("58+3=61", $true, $true, $false, $false, $false, $true, 1, $false, "95-60=35", 2) | Out-Null
$d.Content.Find.Execute("98-82=16", $true, $true, $false, $false, $false, $true, 1, $false, "39-25=14", 2) | Out-Null
$d.Content.Find.Execute("82-68=14", $true, $true, $false, $false, $false, $true, 1, $false, "42-24=18", 2) | Out-Null
$d.Content.Find.Execute("2+39=41", $true, $true, $false, $false, $false, $true, 1, $false, "25-21=4", 2) | Out-Null
$d.Content.Find.Execute("74-12=62", $true, $true, $false, $false, $false, $true, 1, $false, "42+38=80", 2) | Out-Null
$d.Content.Find.Execute("71-14=57", $true, $true, $false, $false, $false, $true, 1, $false, "19+47=66", 2) | Out-Null
$d.Content.Find.Execute("88-79=9", $true, $true, $false, $false, $false, $true, 1, $false, "4+88=92", 2) | Out-Null
$d.Content.Find.Execute("10+2=12", $true, $true, $false, $false, $false, $true, 1, $false, "11-10=1", 2) | Out-Null
$d.Content.Find.Execute("47-29=18", $true, $true, $false, $false, $false, $true, 1, $false, "71-66=5", 2) | Out-Null
$d.Content.Find.Execute("53-27=26", $true, $true, $false, $false, $false, $true, 1, $false, "55-44=11", 2) | Out-Null
$d.Content.Find.Execute("93-62=31", $true, $true, $false, $false, $false, $true, 1, $false, "62-60=2", 2) | Out-Null
$d.Content.Find.Execute("50-11=39", $true, $true, $false, $false, $false, $true, 1, $false, "40-7=33", 2) | Out-Null
$d.Content.Find.Execute("81+13=94", $true, $true, $false, $false, $false, $true, 1, $false, "6+27=33", 2) | Out-Null
$d.Content.Find.Execute("8+47=55", $true, $true, $false, $false, $false, $true, 1, $false, "1+37=38", 2) | Out-Null
